# Publish pass: convert heading paragraphs into Markdown-style "#"/"##"
# prefixed body paragraphs (Heading1 -> FirstParagraph w/ "#" prefix,
# Heading2 -> BodyText w/ "##" prefix, except "Feedback" keeps
# FirstParagraph), and promote the paragraph that used to immediately
# follow each heading from FirstParagraph style to BodyText style.
#
# Paragraph indices are stable (no paragraphs are added or removed),
# so we drive everything off $d.Paragraphs.Item(<n>).

$d = $word.ActiveDocument

# --- Headings: restyle + prepend the Markdown marker -----------------

# 1) "Equipment Reassignment Guide" (Heading 1 -> First Paragraph, "#")
$p = $d.Paragraphs.Item(1)
$p.Range.InsertBefore("#")
$p.Style = "First Paragraph"

# 2) "Introduction" (Heading 2 -> Body Text, "##")
$p = $d.Paragraphs.Item(2)
$p.Range.InsertBefore("##")
$p.Style = "Body Text"

# 4) "Who is this for?" (Heading 2 -> Body Text, "##")
$p = $d.Paragraphs.Item(4)
$p.Range.InsertBefore("##")
$p.Style = "Body Text"

# 10) "Returning Equipment" (Heading 2 -> Body Text, "##")
$p = $d.Paragraphs.Item(10)
$p.Range.InsertBefore("##")
$p.Style = "Body Text"

# 19) "Equipment Reassignment" (Heading 2 -> Body Text, "##")
$p = $d.Paragraphs.Item(19)
$p.Range.InsertBefore("##")
$p.Style = "Body Text"

# 24) "Equipment that cannot be reused" (Heading 2 -> Body Text, "##")
$p = $d.Paragraphs.Item(24)
$p.Range.InsertBefore("##")
$p.Style = "Body Text"

# 27) "Leased equipment" (Heading 2 -> Body Text, "##")
$p = $d.Paragraphs.Item(27)
$p.Range.InsertBefore("##")
$p.Style = "Body Text"

# 29) "Contacts" (Heading 2 -> Body Text, "##")
$p = $d.Paragraphs.Item(29)
$p.Range.InsertBefore("##")
$p.Style = "Body Text"

# 32) "Feedback" (Heading 2 -> First Paragraph, "##")
$p = $d.Paragraphs.Item(32)
$p.Range.InsertBefore("##")
$p.Style = "First Paragraph"

# --- Former lead-in paragraphs: First Paragraph -> Body Text ---------

$d.Paragraphs.Item(3).Style = "Body Text"
$d.Paragraphs.Item(5).Style = "Body Text"
$d.Paragraphs.Item(11).Style = "Body Text"
$d.Paragraphs.Item(20).Style = "Body Text"
$d.Paragraphs.Item(25).Style = "Body Text"
$d.Paragraphs.Item(28).Style = "Body Text"
$d.Paragraphs.Item(30).Style = "Body Text"
$d.Paragraphs.Item(33).Style = "Body Text"

# --- Hyperlink targets stay the same; only relationship ids shift as a
# side effect of the publishing pipeline regenerating the package. That
# is not reachable through the Word object model (no r:id surface), so
# it is intentionally left alone here.
